$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3624.25
$ws.Range("I40").Value = 3248.5
$ws.Range("K40").Value = 3248.5
$ws.Range("M40").Value = -3073.5
$ws.Range("H43").Value = 3377.5557
$ws.Range("I43").Value = 2933
$ws.Range("J43").Value = 3599.8333
$ws.Range("K43").Value = 2933
$ws.Range("L43").Value = 3599.8333
$ws.Range("M43").Value = -2864
$ws.Range("N43").Value = -3737.8333
$ws.Range("H51").Value = 12333.333
$ws.Range("I51").Value = 13000.5
$ws.Range("J51").Value = 10999
$ws.Range("K51").Value = 13000.5
$ws.Range("L51").Value = 10999
$ws.Range("M51").Value = -12516.5
$ws.Range("N51").Value = -11967
$ws.Range("H74").Value = 6999.2
$ws.Range("I74").Value = 6249
$ws.Range("K74").Value = 6249
$ws.Range("M74").Value = -5313
$ws.Range("H76").Value = 4679.8
$ws.Range("J76").Value = 4349.75
$ws.Range("L76").Value = 4349.75
$ws.Range("N76").Value = -4979.75
$ws.Range("H77").Value = 6999.2
$ws.Range("I77").Value = 6249
$ws.Range("K77").Value = 31245
$ws.Range("M77").Value = -26565
$ws.Range("H79").Value = 4679.8
$ws.Range("J79").Value = 4349.75
$ws.Range("L79").Value = 4349.75
$ws.Range("N79").Value = -6533.75
$ws.Range("H112").Value = 3413.2942
$ws.Range("J112").Value = 3413.2942
$ws.Range("L112").Value = 10239.8826
$ws.Range("N112").Value = -12455.8826
$ws.Range("H137").Value = 10948.167
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 12837.8
$ws.Range("K137").Value = 4500
$ws.Range("L137").Value = 38513.39999999999
$ws.Range("M137").Value = -1950
$ws.Range("N137").Value = -43613.39999999999
$ws.Range("H138").Value = 4508.1094
$ws.Range("I138").Value = 2785.8262
$ws.Range("J138").Value = 5300.36
$ws.Range("K138").Value = 8357.4786
$ws.Range("L138").Value = 15901.08
$ws.Range("M138").Value = -3217.4786
$ws.Range("N138").Value = -26181.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23928.98
$ws.Range("I32").Value = 26773.422
$ws.Range("K32").Value = 26773.422
$ws.Range("M32").Value = -26486.422
$ws.Range("H45").Value = 2683.1428
$ws.Range("I45").Value = 857.0909
$ws.Range("J45").Value = 4691.8
$ws.Range("K45").Value = 857.0909
$ws.Range("L45").Value = 4691.8
$ws.Range("M45").Value = -480.0909
$ws.Range("N45").Value = -5445.8
$ws.Range("H61").Value = 2999
$ws.Range("I61").Value = 1998
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1998
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1786
$ws.Range("N61").Value = -4424
$ws.Range("H74").Value = 126693.875
$ws.Range("I74").Value = 201509.6
$ws.Range("J74").Value = 2001
$ws.Range("K74").Value = 201509.6
$ws.Range("L74").Value = 2001
$ws.Range("M74").Value = -200635.6
$ws.Range("N74").Value = -3749
$ws.Range("H77").Value = 126693.875
$ws.Range("I77").Value = 201509.6
$ws.Range("J77").Value = 2001
$ws.Range("K77").Value = 1007548
$ws.Range("L77").Value = 10005
$ws.Range("M77").Value = -1003180
$ws.Range("N77").Value = -18741
$ws.Range("H132").Value = 66117.31
$ws.Range("I132").Value = 85624.25
$ws.Range("J132").Value = 7596.5
$ws.Range("K132").Value = 256872.75
$ws.Range("L132").Value = 22789.5
$ws.Range("M132").Value = -254342.75
$ws.Range("N132").Value = -27849.5
$ws.Range("H136").Value = 2999
$ws.Range("I136").Value = 1998
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 5994
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -3444
$ws.Range("N136").Value = -17100
$ws.Range("H138").Value = 200000
$ws.Range("J138").Value = 200000
$ws.Range("L138").Value = 200000
$ws.Range("N138").Value = -210280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 65
$ws.Range("I7").Value = 47.1
$ws.Range("K7").Value = 47.1
$ws.Range("M7").Value = 65.90000000000001
$ws.Range("H99").Value = 4559.2
$ws.Range("J99").Value = 3500
$ws.Range("L99").Value = 3500
$ws.Range("N99").Value = -6496
$ws.Range("H122").Value = 2327.5386
$ws.Range("I122").Value = 2125.4
$ws.Range("K122").Value = 6376.200000000001
$ws.Range("M122").Value = -3926.200000000001
$ws.Range("H126").Value = 4559.2
$ws.Range("J126").Value = 3500
$ws.Range("L126").Value = 10500
$ws.Range("N126").Value = -15440
$ws.Range("H132").Value = 1613.3214
$ws.Range("I132").Value = 1444.5471
$ws.Range("K132").Value = 4333.6413
$ws.Range("M132").Value = -1803.6413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100.53333
$ws.Range("I2").Value = 137.6
$ws.Range("J2").Value = 26.4
$ws.Range("K2").Value = 825.5999999999999
$ws.Range("L2").Value = 158.4
$ws.Range("M2").Value = -712.5999999999999
$ws.Range("N2").Value = -384.4
$ws.Range("H4").Value = 761270.1
$ws.Range("J4").Value = 1148.5
$ws.Range("L4").Value = 3445.5
$ws.Range("N4").Value = -3669.5
$ws.Range("H37").Value = 68328.57000000001
$ws.Range("J37").Value = 68328.57000000001
$ws.Range("L37").Value = 204985.71
$ws.Range("N37").Value = -205209.71
$ws.Range("H55").Value = 6928.143
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 7249.5
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 21748.5
$ws.Range("M55").Value = -14823
$ws.Range("N55").Value = -22102.5
$ws.Range("H94").Value = 14512.75
$ws.Range("I94").Value = 12012
$ws.Range("J94").Value = 15346.333
$ws.Range("K94").Value = 36036
$ws.Range("L94").Value = 46038.999
$ws.Range("M94").Value = -35360
$ws.Range("N94").Value = -47390.999
$ws.Range("H122").Value = 405.6875
$ws.Range("I122").Value = 361.77777
$ws.Range("J122").Value = 462.14285
$ws.Range("K122").Value = 3255.99993
$ws.Range("L122").Value = 4159.28565
$ws.Range("M122").Value = -805.9999299999999
$ws.Range("N122").Value = -9059.28565
$ws.Range("H131").Value = 2228047.5
$ws.Range("I131").Value = 2279.3635
$ws.Range("J131").Value = 2948149
$ws.Range("K131").Value = 6838.0905
$ws.Range("L131").Value = 8844447
$ws.Range("M131").Value = -1798.0905
$ws.Range("N131").Value = -8854527

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3094.5
$ws.Range("I80").Value = 2931.2856
$ws.Range("J80").Value = 3221.4443
$ws.Range("K80").Value = 2931.2856
$ws.Range("L80").Value = 3221.4443
$ws.Range("M80").Value = -1933.2856
$ws.Range("N80").Value = -5217.4443
$ws.Range("H83").Value = 3094.5
$ws.Range("I83").Value = 2931.2856
$ws.Range("J83").Value = 3221.4443
$ws.Range("K83").Value = 14656.428
$ws.Range("L83").Value = 16107.2215
$ws.Range("M83").Value = -9664.428
$ws.Range("N83").Value = -26091.2215
$ws.Range("H102").Value = 3917.4
$ws.Range("I102").Value = 2396.75
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 2396.75
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -774.75
$ws.Range("N102").Value = -13244
$ws.Range("H126").Value = 7104.5713
$ws.Range("I126").Value = 3946.4
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 11839.2
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = -9369.200000000001
$ws.Range("N126").Value = -49940
$ws.Range("H132").Value = 106498.7
$ws.Range("I132").Value = 145752.58
$ws.Range("K132").Value = 437257.74
$ws.Range("M132").Value = -434727.74
$ws.Range("H141").Value = 69000
$ws.Range("J141").Value = 69000
$ws.Range("L141").Value = 69000
$ws.Range("N141").Value = -79360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5517.6665
$ws.Range("J40").Value = 8200
$ws.Range("L40").Value = 8200
$ws.Range("N40").Value = -8472
$ws.Range("H42").Value = 30122
$ws.Range("J42").Value = 32499
$ws.Range("L42").Value = 32499
$ws.Range("N42").Value = -33625
$ws.Range("H46").Value = 14744.111
$ws.Range("J46").Value = 5813.857
$ws.Range("L46").Value = 5813.857
$ws.Range("N46").Value = -6189.857
$ws.Range("H49").Value = 30122
$ws.Range("J49").Value = 32499
$ws.Range("L49").Value = 32499
$ws.Range("N49").Value = -32793
$ws.Range("H61").Value = 2007.6086
$ws.Range("I61").Value = 2007.6086
$ws.Range("K61").Value = 2007.6086
$ws.Range("M61").Value = -1805.6086
$ws.Range("H68").Value = 6933
$ws.Range("J68").Value = 6933
$ws.Range("L68").Value = 6933
$ws.Range("N68").Value = -8431
$ws.Range("H71").Value = 6933
$ws.Range("J71").Value = 6933
$ws.Range("L71").Value = 34665
$ws.Range("N71").Value = -42153
$ws.Range("H113").Value = 2007.6086
$ws.Range("I113").Value = 2007.6086
$ws.Range("K113").Value = 2007.6086
$ws.Range("M113").Value = 162.3914
$ws.Range("H132").Value = 67223.58
$ws.Range("I132").Value = 96134.53999999999
$ws.Range("J132").Value = 4583.1665
$ws.Range("K132").Value = 288403.62
$ws.Range("L132").Value = 13749.4995
$ws.Range("M132").Value = -285873.62
$ws.Range("N132").Value = -18809.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 443285
$ws.Range("I126").Value = 443285
$ws.Range("K126").Value = 1329855
$ws.Range("M126").Value = -1327385
$ws.Range("H132").Value = 75929.39
$ws.Range("I132").Value = 79086.71000000001
$ws.Range("K132").Value = 237260.13
$ws.Range("M132").Value = -234730.13
